$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 520, pushing existing rows 520:614 down to 521:615
$ws.Rows("520:520").Insert()

# Populate the newly inserted row 520 with the new record
$ws.Range("A520").Value = 8
$ws.Range("B520").Value = "Terminal La Palmera de La Serena"
$ws.Range("C520").Value = "Coquimbo"
$ws.Range("D520").Value = 45015
$ws.Range("E520").Value = 4
$ws.Range("F520").Value = 100114001
$ws.Range("G520").Value = "Papa"
$ws.Range("H520").Value = "Asterix"
$ws.Range("I520").Value = "1a (cosecha)"
$ws.Range("J520").Value = 2000
$ws.Range("K520").Value = 11000
$ws.Range("L520").Value = 12000
$ws.Range("M520").Value = 11500
$ws.Range("N520").Value = "$/saco 25 kilos"
$ws.Range("O520").Value = "Provincia de Melipilla"
$ws.Range("P520").Value = 460
$ws.Range("Q520").Value = 25
$ws.Range("R520").Value = "Hortaliza"
